$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Freddy Velez"
$ws.Range("A2").Value = "4FD889D140"
$ws.Range("A3").Value = "IN -> 2017/02/13 18:45"
